$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "extr*" contingency rows (originally rows 8-15) need to move down by two
# rows to make room for two new "line7"/"line8" rows that get inserted right
# after the existing "line6" row. Shift bottom-up (15 -> 17, 8 -> 10) using
# Value2 to read (Value's getter is unreliable for reads in this host) so we
# never clobber data we still need to read.
# ---------------------------------------------------------------------------
for ($r = 15; $r -ge 8; $r--) {
    $dst = $r + 2
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dst, 4).Value = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($dst, 5).Value = $ws.Cells.Item($r, 5).Value2
}

# The two rows that used to hold row 15 and row 14 (now rows 17 and 16) are
# brand new to the sheet's dimension, so they come out with no formatting.
# Pull the cell style from the existing "A" column (row 2) which already
# carries the correct style used throughout that column.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Fill in the two new "line" rows (8 and 9).
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $false

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# ---------------------------------------------------------------------------
# Re-index column A (the sequence number) and refresh C/D/E for the shifted
# "extr*" rows (10-17) to match the final, "rene fine"-refined contingency
# numbers.
# ---------------------------------------------------------------------------
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $false

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $true

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $true

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $true

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false
